# The deck currently uses the "Red Violet" variant of the "Integral" theme
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) on its slide master/theme.
# The target revision swaps the live theme's colour scheme back to the
# stock Office colour scheme (the colours that the deck's other, inactive
# theme part already carries), leaving the font scheme / format scheme
# (both already "Office") untouched.
#
# PowerPoint's object model exposes the 12 theme colours as an indexed
# collection (1=dk1, 2=lt1, 3=dk2, 4=lt2, 5-10=accent1-6, 11=hlink,
# 12=folHlink) via ThemeColorScheme.Colors(i).RGB. Driving that collection
# is the supported, scriptable equivalent of picking the built-in "Office"
# colour variant from the Design > Variants > Colors gallery.

$p = $ppt.ActivePresentation

# RGB() packs as 0x00BBGGRR, matching the OLE colour convention used by
# ThemeColor.RGB / RGBColor.RGB in the PowerPoint object model.
function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office" theme colours (hex -> decimal RGB() value):
#   dk1      000000
#   lt1      FFFFFF
#   dk2      44546A
#   lt2      E7E6E6
#   accent1  5B9BD5
#   accent2  ED7D31
#   accent3  A5A5A5
#   accent4  FFC000
#   accent5  4472C4
#   accent6  70AD47
#   hlink    0563C1
#   folHlink 954F72
$officeColors = @(
    (RGB 0x00 0x00 0x00),  # 1  dk1
    (RGB 0xFF 0xFF 0xFF),  # 2  lt1
    (RGB 0x44 0x54 0x6A),  # 3  dk2
    (RGB 0xE7 0xE6 0xE6),  # 4  lt2
    (RGB 0x5B 0x9B 0xD5),  # 5  accent1
    (RGB 0xED 0x7D 0x31),  # 6  accent2
    (RGB 0xA5 0xA5 0xA5),  # 7  accent3
    (RGB 0xFF 0xC0 0x00),  # 8  accent4
    (RGB 0x44 0x72 0xC4),  # 9  accent5
    (RGB 0x70 0xAD 0x47),  # 10 accent6
    (RGB 0x05 0x63 0xC1),  # 11 hlink
    (RGB 0x95 0x4F 0x72)   # 12 folHlink
)

# Apply the new colour scheme to the presentation's theme (ThemeColorScheme
# is reachable from any slide and edits the one shared theme used across
# the deck, including the slide master).
$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
